$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "31/01/2022"
$newRow.Cells.Item(2).Range.Text = "1 Hour 30 Minutes"
$newRow.Cells.Item(3).Range.Text = "Simulation – Objective 3"
$newRow.Cells.Item(4).Range.Text = "Added the ability for provinces to spawn empires and add land to their empire using console commands. This is largely untested thusfar and does not save as of yet."
